# The original workbook had several shared-string cell values containing
# embedded line breaks (entered in Excel with Alt+Enter). This edit removes
# those line breaks, replacing them with a single space, so each of the
# following values becomes one line of text instead of two or three.
#
# "Pediatric Vaccine" sheet:
#   D8: "5 pack - 1 dose\nT-L syringes. No Needle"
#       -> "5 pack - 1 dose T-L syringes. No Needle"
#
# "Pediatric Influenza Vaccine" sheet:
#   B3:  "Fluzone\nPediatric dose\nNo Preservative"
#        -> "Fluzone Pediatric dose No Preservative"
#   B6:  "Fluarix\nPreservative-Free"      -> "Fluarix Preservative-Free"
#   B9:  "FluMist\nNo Preservative"        -> "FluMist No Preservative"
#   B10: "Afluria\nNo Preservative"        -> "Afluria No Preservative"
#   H10: "Merck\n(CSL product)"            -> "Merck (CSL product)"
#   H11: "Merck\n(CSL product)"            -> "Merck (CSL product)"
#   B12: "Afluria\nNo Preservative"        -> "Afluria No Preservative"
#   H12: "Merck\n(CSL product)"            -> "Merck (CSL product)"
#
# "Adult Influenza Vaccine" sheet:
#   B5:  "Agriflu\nNo Preservative"        -> "Agriflu No Preservative"
#   B7:  "Fluvirin\nPreservative-free"     -> "Fluvirin Preservative-free"
#   B8:  "Fluarix\nPreservative-free"      -> "Fluarix Preservative-free"
#   B10: "Flumist\nNo Preservative"        -> "Flumist No Preservative"

$wb = $excel.ActiveWorkbook

$wsPed = $wb.Worksheets.Item("Pediatric Vaccine ")
$wsPed.Range("D8").Value = "5 pack - 1 dose T-L syringes. No Needle"

$wsPedFlu = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$wsPedFlu.Range("B3").Value = "Fluzone Pediatric dose No Preservative"
$wsPedFlu.Range("B6").Value = "Fluarix Preservative-Free"
$wsPedFlu.Range("B9").Value = "FluMist No Preservative"
$wsPedFlu.Range("B10").Value = "Afluria No Preservative"
$wsPedFlu.Range("H10").Value = "Merck (CSL product)"
$wsPedFlu.Range("H11").Value = "Merck (CSL product)"
$wsPedFlu.Range("B12").Value = "Afluria No Preservative"
$wsPedFlu.Range("H12").Value = "Merck (CSL product)"

$wsAdultFlu = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$wsAdultFlu.Range("B5").Value = "Agriflu No Preservative"
$wsAdultFlu.Range("B7").Value = "Fluvirin Preservative-free"
$wsAdultFlu.Range("B8").Value = "Fluarix Preservative-free"
$wsAdultFlu.Range("B10").Value = "Flumist No Preservative"

Write-Output "done"
